$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Staging area (far outside used range) for numeric-looking text values so they are not
# auto-coerced to numbers when written; a formula yielding a text result, copied and
# pasted as values into the destination, preserves text type without touching styles.
$stageRow = 1000

# ---- row 12 ----
$ws.Cells.Item(12, 1).Value = "3RWWKEG0@testmail.com"
$ws.Cells.Item(12, 2).Value = "Daniel"
$ws.Cells.Item(12, 3).Value = "apellidos"
$ws.Cells.Item($stageRow, 1).Formula = "=""84851"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1903615730"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 6).Value = "Cll 80 # 106 - 56"
$ws.Cells.Item(12, 7).Value = "cartagena"
$ws.Cells.Item($stageRow, 1).Formula = "=""22883"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(12, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 13 ----
$ws.Cells.Item(13, 1).Value = "PIM9N83O@testmail.com"
$ws.Cells.Item(13, 2).Value = "Daniel"
$ws.Cells.Item(13, 3).Value = "apellidos"
$ws.Cells.Item($stageRow, 1).Formula = "=""20153"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1197460105"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 6).Value = "Cra 90#17-90"
$ws.Cells.Item(13, 7).Value = "rio de janeiro"
$ws.Cells.Item($stageRow, 1).Formula = "=""36854"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(13, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 14 ----
$ws.Cells.Item(14, 1).Value = "CBFSDOJO@testmail.com"
$ws.Cells.Item(14, 2).Value = "Daniel"
$ws.Cells.Item(14, 3).Value = "Lopez"
$ws.Cells.Item($stageRow, 1).Formula = "=""73164"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1306698380"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 6).Value = "Cra 90#17-90"
$ws.Cells.Item(14, 7).Value = "rio de janeiro"
$ws.Cells.Item($stageRow, 1).Formula = "=""90535"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(14, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 15 ----
$ws.Cells.Item(15, 1).Value = "DUQL50SX@testmail.com"
$ws.Cells.Item(15, 2).Value = "Pedro"
$ws.Cells.Item(15, 3).Value = "Lopez"
$ws.Cells.Item($stageRow, 1).Formula = "=""71620"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1935491711"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 6).Value = "Cra 90#17-90"
$ws.Cells.Item(15, 7).Value = "rio de janeiro"
$ws.Cells.Item($stageRow, 1).Formula = "=""20112"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(15, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 16 ----
$ws.Cells.Item(16, 1).Value = "0C6Y6RQN@testmail.com"
$ws.Cells.Item(16, 2).Value = "Gerardo"
$ws.Cells.Item(16, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""11484"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1630656925"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$ws.Cells.Item(16, 6).Value = "cra # 123A - 45"
$ws.Cells.Item(16, 7).Value = "santiago"
$ws.Cells.Item($stageRow, 1).Formula = "=""69919"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(16, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 17 ----
$ws.Cells.Item(17, 1).Value = "334DWD7E@testmail.com"
$ws.Cells.Item(17, 2).Value = "Daniel"
$ws.Cells.Item(17, 3).Value = "Ramos"
$ws.Cells.Item($stageRow, 1).Formula = "=""33268"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1259877314"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 6).Value = "Cra 7 #72-21"
$ws.Cells.Item(17, 7).Value = "rio de janeiro"
$ws.Cells.Item($stageRow, 1).Formula = "=""41892"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(17, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 18 ----
$ws.Cells.Item(18, 1).Value = "YHXTCDP1@testmail.com"
$ws.Cells.Item(18, 2).Value = "Hombres"
$ws.Cells.Item(18, 3).Value = "Ramos"
$ws.Cells.Item($stageRow, 1).Formula = "=""14862"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1971059754"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 6).Value = "Cra 7 #72-21"
$ws.Cells.Item(18, 7).Value = "cali"
$ws.Cells.Item($stageRow, 1).Formula = "=""20413"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(18, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 19 ----
$ws.Cells.Item(19, 1).Value = "HHJ7TGZ4@testmail.com"
$ws.Cells.Item(19, 2).Value = "Daniel"
$ws.Cells.Item(19, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""79629"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1831672394"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 6).Value = "Av 3 # 45 -2"
$ws.Cells.Item(19, 7).Value = "cartagena"
$ws.Cells.Item($stageRow, 1).Formula = "=""64015"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(19, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 20 ----
$ws.Cells.Item(20, 1).Value = "XL2UV568@testmail.com"
$ws.Cells.Item(20, 2).Value = "Daniel"
$ws.Cells.Item(20, 3).Value = "Ramos"
$ws.Cells.Item($stageRow, 1).Formula = "=""46060"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1355870395"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 6).Value = "Cra 7 #72-21"
$ws.Cells.Item(20, 7).Value = "city "
$ws.Cells.Item($stageRow, 1).Formula = "=""90809"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(20, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 21 ----
$ws.Cells.Item(21, 1).Value = "8IAGAA86@testmail.com"
$ws.Cells.Item(21, 2).Value = "Pablo"
$ws.Cells.Item(21, 3).Value = "Gonzalez"
$ws.Cells.Item($stageRow, 1).Formula = "=""48588"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1931774627"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 6).Value = "Cll 90 bis #32 - 43"
$ws.Cells.Item(21, 7).Value = "medellin"
$ws.Cells.Item($stageRow, 1).Formula = "=""27988"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(21, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 22 ----
$ws.Cells.Item(22, 1).Value = "GIQKVG9Y@testmail.com"
$ws.Cells.Item(22, 2).Value = "Sergio"
$ws.Cells.Item(22, 3).Value = "Gomez"
$ws.Cells.Item($stageRow, 1).Formula = "=""50836"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1159795927"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 6).Value = "Cra 90#17-90"
$ws.Cells.Item(22, 7).Value = "pereira"
$ws.Cells.Item($stageRow, 1).Formula = "=""83623"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(22, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 23 ----
$ws.Cells.Item(23, 1).Value = "1LYDZ3F7@testmail.com"
$ws.Cells.Item(23, 2).Value = "Juan"
$ws.Cells.Item(23, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""46128"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1187334222"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$ws.Cells.Item(23, 6).Value = "Cll 80 # 106 - 56"
$ws.Cells.Item(23, 7).Value = "santiago"
$ws.Cells.Item($stageRow, 1).Formula = "=""81809"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(23, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 24 ----
$ws.Cells.Item(24, 1).Value = "TJ30QYYP@testmail.com"
$ws.Cells.Item(24, 2).Value = "Alfredo"
$ws.Cells.Item(24, 3).Value = "Vargas"
$ws.Cells.Item($stageRow, 1).Formula = "=""33602"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1772766936"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$ws.Cells.Item(24, 6).Value = "Cra 80 # 47 -90"
$ws.Cells.Item(24, 7).Value = "pereira"
$ws.Cells.Item($stageRow, 1).Formula = "=""21264"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(24, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 25 ----
$ws.Cells.Item(25, 1).Value = "3Z2PJN9L@testmail.com"
$ws.Cells.Item(25, 2).Value = "Daniel"
$ws.Cells.Item(25, 3).Value = "Gomez"
$ws.Cells.Item($stageRow, 1).Formula = "=""71100"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1699025943"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 6).Value = "Cra 90#17-90"
$ws.Cells.Item(25, 7).Value = "pereira"
$ws.Cells.Item($stageRow, 1).Formula = "=""50222"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(25, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 26 ----
$ws.Cells.Item(26, 1).Value = "JUTK8YMV@testmail.com"
$ws.Cells.Item(26, 2).Value = "Daniel"
$ws.Cells.Item(26, 3).Value = "Rey"
$ws.Cells.Item($stageRow, 1).Formula = "=""52765"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1215577959"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$ws.Cells.Item(26, 6).Value = "Cra 7 #72-21"
$ws.Cells.Item(26, 7).Value = "bogota"
$ws.Cells.Item($stageRow, 1).Formula = "=""83210"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(26, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 27 ----
$ws.Cells.Item(27, 1).Value = "7IXYJFHE@testmail.com"
$ws.Cells.Item(27, 2).Value = "Sergio"
$ws.Cells.Item(27, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""33870"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1534758411"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$ws.Cells.Item(27, 6).Value = "Cra 80 # 47 -90"
$ws.Cells.Item(27, 7).Value = "bogota"
$ws.Cells.Item($stageRow, 1).Formula = "=""21320"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(27, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 28 ----
$ws.Cells.Item(28, 1).Value = "3H1OACDQ@testmail.com"
$ws.Cells.Item(28, 2).Value = "Juan"
$ws.Cells.Item(28, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""57095"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1525507756"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 6).Value = "Cll 90 bis #32 - 43"
$ws.Cells.Item(28, 7).Value = "medellin"
$ws.Cells.Item($stageRow, 1).Formula = "=""71145"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(28, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 29 ----
$ws.Cells.Item(29, 1).Value = "0M82CY92@testmail.com"
$ws.Cells.Item(29, 2).Value = "Hombres"
$ws.Cells.Item(29, 3).Value = "Penagos"
$ws.Cells.Item($stageRow, 1).Formula = "=""60307"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1259630638"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$ws.Cells.Item(29, 6).Value = "cra # 123A - 45"
$ws.Cells.Item(29, 7).Value = "santiago"
$ws.Cells.Item($stageRow, 1).Formula = "=""17172"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(29, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 30 ----
$ws.Cells.Item(30, 1).Value = "WKGFKFLD@testmail.com"
$ws.Cells.Item(30, 2).Value = "Daniel"
$ws.Cells.Item(30, 3).Value = "Vargas"
$ws.Cells.Item($stageRow, 1).Formula = "=""99441"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1038202514"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 6).Value = "direccion"
$ws.Cells.Item(30, 7).Value = "cartagena"
$ws.Cells.Item($stageRow, 1).Formula = "=""20427"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(30, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 31 ----
$ws.Cells.Item(31, 1).Value = "AT8G3AZT@testmail.com"
$ws.Cells.Item(31, 2).Value = "Pedro"
$ws.Cells.Item(31, 3).Value = "Gonzalez"
$ws.Cells.Item($stageRow, 1).Formula = "=""35358"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1781465293"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$ws.Cells.Item(31, 6).Value = "cra # 123A - 45"
$ws.Cells.Item(31, 7).Value = "santiago"
$ws.Cells.Item($stageRow, 1).Formula = "=""34980"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(31, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 32 ----
$ws.Cells.Item(32, 1).Value = "1598IL61@testmail.com"
$ws.Cells.Item(32, 2).Value = "Sergio"
$ws.Cells.Item(32, 3).Value = "Gomez"
$ws.Cells.Item($stageRow, 1).Formula = "=""18960"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1694083677"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$ws.Cells.Item(32, 6).Value = "direccion"
$ws.Cells.Item(32, 7).Value = "city "
$ws.Cells.Item($stageRow, 1).Formula = "=""32054"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(32, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 33 ----
$ws.Cells.Item(33, 1).Value = "OZJGL94W@testmail.com"
$ws.Cells.Item(33, 2).Value = "Alfredo"
$ws.Cells.Item(33, 3).Value = "Sanchez"
$ws.Cells.Item($stageRow, 1).Formula = "=""51168"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1887775936"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 6).Value = "Cra 7 #72-21"
$ws.Cells.Item(33, 7).Value = "pereira"
$ws.Cells.Item($stageRow, 1).Formula = "=""11613"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(33, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

# ---- row 34 ----
$ws.Cells.Item(34, 1).Value = "AV3YPEC2@testmail.com"
$ws.Cells.Item(34, 2).Value = "Daniel"
$ws.Cells.Item(34, 3).Value = "Gonzalez"
$ws.Cells.Item($stageRow, 1).Formula = "=""61757"""
$ws.Cells.Item($stageRow, 2).Formula = "=""1383261546"""
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 6).Value = "cra # 123A - 45"
$ws.Cells.Item(34, 7).Value = "bogota"
$ws.Cells.Item($stageRow, 1).Formula = "=""78562"""
$ws.Cells.Item($stageRow, 1).Copy()
$ws.Cells.Item(34, 8).PasteSpecial(-4163)
$ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, 2)).Clear()
$excel.CutCopyMode = $false

